$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "43.696.58"
Set-TextCell "E2" "  +2.52%  "
Set-TextCell "D3" "2.201.44"
Set-TextCell "E3" "  +0.20%  "
Set-TextCell "E4" "  -0.03%  "
Set-TextCell "D5" "258.04"
Set-TextCell "E5" "  +2.27%  "
Set-TextCell "D6" "84.31"
Set-TextCell "E6" "  +11.78%  "
Set-TextCell "D7" "0.615"
Set-TextCell "E7" "  +0.01%  "
Set-TextCell "E8" "  -0.12%  "
Set-TextCell "D9" "0.599"
Set-TextCell "E9" "  +2.28%  "
Set-TextCell "D10" "44.43"
Set-TextCell "E10" "  +9.35%  "
Set-TextCell "E11" "  +0.53%  "
Set-TextCell "D12" "7.22"
Set-TextCell "E12" "  +5.40%  "
Set-TextCell "E13" "  +2.47%  "
Set-TextCell "D14" "2.530.25"
Set-TextCell "E14" "  +0.43%  "
Set-TextCell "D15" "14.31"
Set-TextCell "E15" "  +0.47%  "
Set-TextCell "D16" "2.202.61"
Set-TextCell "E16" "  +0.95%  "
Set-TextCell "D17" "0.780"
Set-TextCell "E17" "  +0.69%  "
Set-TextCell "D18" "43.637.99"
Set-TextCell "E18" "  +2.57%  "
Set-TextCell "E19" "  +0.85%  "
Set-TextCell "D20" "69.73"
Set-TextCell "E20" "  -1.65%  "
Set-TextCell "D21" "5.92"
Set-TextCell "E21" "  +0.18%  "
Set-TextCell "D22" "2.38"
Set-TextCell "E22" "  +9.96%  "
Set-TextCell "D23" "231.81"
Set-TextCell "E23" "  +1.78%  "
Set-TextCell "D24" "9.05"
Set-TextCell "E24" "  -4.86%  "
Set-TextCell "E25" "  -0.09%  "
Set-TextCell "D26" "3.52"
Set-TextCell "E26" "  +4.22%  "
Set-TextCell "D27" "10.67"
Set-TextCell "E27" "  +1.64%  "
Set-TextCell "D28" "39.09"
Set-TextCell "E28" "  +0.88%  "
Set-TextCell "E29" "  +2.98%  "
Set-TextCell "D31" "173.77"
Set-TextCell "E31" "  +0.32%  "
Set-TextCell "D32" "20.40"
Set-TextCell "E32" "  +1.34%  "
Set-TextCell "D33" "0.0858"
Set-TextCell "E33" "  +4.39%  "
Set-TextCell "D34" "5.32"
Set-TextCell "E34" "  +2.94%  "
Set-TextCell "E35" "  +1.57%  "
Set-TextCell "E36" "  +2.96%  "
Set-TextCell "E37" "  +6.32%  "
Set-TextCell "D38" "4.49"
Set-TextCell "E38" "  +6.01%  "
Set-TextCell "D39" "12.52"
Set-TextCell "E39" "  +2.55%  "
Set-TextCell "D40" "2.87"
Set-TextCell "E40" "  +9.25%  "
Set-TextCell "E41" "  +0.93%  "
Set-TextCell "D42" "63.12"
Set-TextCell "E42" "  +5.95%  "
Set-TextCell "D43" "5.49"
Set-TextCell "E43" "  +5.76%  "
Set-TextCell "D44" "0.200"
Set-TextCell "E44" "  +2.35%  "
Set-TextCell "E45" "  +1.22%  "
Set-TextCell "D46" "0.0979"
Set-TextCell "E46" "  +0.04%  "
Set-TextCell "D47" "99.84"
Set-TextCell "E47" "  -2.21%  "
Set-TextCell "E48" "  +5.91%  "
Set-TextCell "E49" "  +1.35%  "
Set-TextCell "D50" "0.436"
Set-TextCell "E50" "  -5.47%  "
Set-TextCell "D51" "1.50"
Set-TextCell "E51" "  +9.08%  "
